$d = $word.ActiveDocument

# This tc_p016r.docx transcription renders inline "rendition spec" markup
# (<tl>...</tl>, <m>...</m>, <pa>...</pa>, <lb/>, ...) as literal,
# specially-coloured/"Courier New" text runs right next to the plain
# transcription text they annotate. The commit removes the <tl>/</tl> and
# <m>/</m> wrapper runs around three occurrences of "fonte"/"fontes",
# leaving the plain word itself (and its own run formatting) untouched.
#
# Remove-Wrapper finds a uniquely-identifying chunk of text containing
# "<openTag>content</closeTag>", then deletes just the open/close tag
# runs (back to front, so offsets don't shift under us) while leaving
# the "content" text completely alone. When the runs adjacent to the
# deleted tag share identical formatting, Word naturally coalesces them
# into a single run - matching the target diff exactly.
function Remove-Wrapper($searchText, $openTag, $content, $closeTag) {
    $d = $word.ActiveDocument
    $r = $d.Content
    $found = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Remove-Wrapper: could not find '$searchText'"
    }

    $openRange = $r.Duplicate
    $openRange.End = $openRange.Start + $openTag.Length

    $closeRange = $r.Duplicate
    $closeRange.Start = $closeRange.Start + $openTag.Length + $content.Length
    $closeRange.End = $closeRange.Start + $closeTag.Length

    # Delete the closing tag first so the opening tag's offsets stay valid.
    $closeRange.Delete()
    $openRange.Delete()
}

# "bouche de leur <tl>fonte</tl> deulx ou" -> "bouche de leur fonte deulx ou"
Remove-Wrapper "<tl>fonte</tl> deulx ou" "<tl>" "fonte" "</tl>"

# "faire courre la <m>fonte</m>" -> "faire courre la fonte"
Remove-Wrapper "<m>fonte</m>" "<m>" "fonte" "</m>"

# "e aulx aultres <tl>fontes</tl> Affin de mectre" -> "e aulx aultres fontes Affin de mectre"
Remove-Wrapper "<tl>fontes</tl> Affin de mectre" "<tl>" "fontes" "</tl>"
